# Deploying to gh-pages from @ d0bca759f32cb3bb315e1f39e6e0ba8a25bbfc1f 🚀
# Adds the 2022 data column (S) to the Hepatitis B incidence sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at S. Excel copies formatting from the column
# immediately to the left (R), which already carries the correct
# number format / borders for every row in this table.
$ws.Columns("S").Insert()

# Header (year) cell.
$ws.Range("S3").Value2 = 2022

# Data values for 2022, one per region/category row.
$ws.Range("S4").Value2 = 1.9210869108320343
$ws.Range("S5").Value2 = 1.020872301352429
$ws.Range("S6").Value2 = 2.8415499553180767
$ws.Range("S7").Value2 = 1.5924017665043597
$ws.Range("S8").Value2 = 2.5011433798307796
$ws.Range("S9").Value2 = 0.70098698968147144
$ws.Range("S10").Value2 = 2.2312343573160249
$ws.Range("S11").Value2 = 2.4764236727529938
$ws.Range("S12").Value2 = 1.9888745417939038
$ws.Range("S13").Value2 = 1.3057776932131271
$ws.Range("S14").Value2 = 2.6056788910230639
$ws.Range("S15").Value2 = 0
$ws.Range("S16").Value2 = 0.65058422463372112
$ws.Range("S17").Value2 = 0.65686622262510019
$ws.Range("S18").Value2 = 0.64442124527961442
$ws.Range("S19").Value2 = 2.5553368555544047
$ws.Range("S20").Value2 = 1.807815324711445
$ws.Range("S21").Value2 = 3.2928586128833093
$ws.Range("S22").Value2 = 1.8387963974300983
$ws.Range("S23").Value2 = 2.2260807622100529
$ws.Range("S24").Value2 = 1.4582467499325562
$ws.Range("S25").Value2 = 1.2245886088767601
$ws.Range("S26").Value2 = 1.3105423773238725
$ws.Range("S27").Value2 = 1.1375464261135158
$ws.Range("S28").Value2 = 2.4791112740241377
$ws.Range("S29").Value2 = 2.4279584268771761
$ws.Range("S30").Value2 = 2.5408788313520994
$ws.Range("S31").Value2 = 1.1238322680339958
$ws.Range("S32").Value2 = 0.57553956834532372
$ws.Range("S33").Value2 = 1.6467682173734046

# The "region total" rows (7, 10, 13, 16, 19, 22, 25, 28, 31) use a bold
# variant of the numeric style elsewhere in this column, so give their
# new 2022 cells the same bold treatment (creates a dedicated style,
# distinct from the regular data-row style).
$boldRows = 7, 10, 13, 16, 19, 22, 25, 28, 31
foreach ($r in $boldRows) {
    $cell = $ws.Range("S$r")
    $cell.Font.Bold = $true
    $cell.WrapText = $false
}

# Move the active selection, matching the author's post-edit cursor spot.
$ws.Range("T3").Select()

Write-Output "Added 2022 column (S) to sheet $($ws.Name)"
